# ESM_input_test.xlsx update
#
# 1) "Views" sheet: add a new row ("flow_cy" / "flow_cyt") documenting the
#    new flow-cytometry transformation/view, and leave the selection back at A2.
# 2) "ID" sheet: rename the header columns from "Current ID"/"Replaced ID"
#    to "Current"/"Target", leaving the selection on B2.
# 3) "Transformations" sheet: add the new "flow_cyt" transformation row that
#    calls process_fcs(...) on plate_02. This is the last sheet touched, so
#    it ends up the active tab, with the selection sitting on B18 (just below
#    the newly entered row).

$wb = $excel.ActiveWorkbook

# --- Views sheet -----------------------------------------------------------
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Range("A4").Value = "flow_cy"
$wsViews.Range("B4").Value = "flow_cyt"
[void]$wsViews.Range("A2").Select()

# --- ID sheet ---------------------------------------------------------------
$wsID = $wb.Worksheets.Item("ID")
$wsID.Range("A1").Value = "Current"
$wsID.Range("B1").Value = "Target"
[void]$wsID.Range("B2").Select()

# --- Transformations sheet (ends as the active sheet) -----------------------
$wsTransformations = $wb.Worksheets.Item("Transformations")
$wsTransformations.Activate()
$wsTransformations.Range("A17").Value = "flow_cyt"
$wsTransformations.Range("B17").Value = 'process_fcs("plate_02",["FSC-A","SSC-A"],["BL1-H"])'
[void]$wsTransformations.Range("B18").Select()
